# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") previously stored a value derived from
# the old "Strike#" computation. It has been regenerated so that it now holds
# the newly calculated K values (s_vals) for each game row (rows 2-56).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..56 (in row order)
$kValues = @(2,3,2,3,4,0,0,2,1,1,0,2,0,0,1,1,2,1,2,0,0,0,3,1,2,1,3,0,0,1,0,1,1,0,0,2,0,2,4,2,1,0,2,2,1,0,4,1,1,2,3,2,1,0,2)

$startRow = 2
$endRow = 56
$rowCount = $endRow - $startRow + 1

$arr = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i, 0] = $kValues[$i]
}

$range = $ws.Range("G$startRow`:G$endRow")
$range.Value = $arr

Write-Host "Updated K (column G) values for rows $startRow-$endRow"
